$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D; this shifts existing columns D:K to F:M,
# matching the quarter-over-quarter roll seen in the diff (two new reporting
# periods added: 2018-12-31 and 2018-09-30, i.e. serials 43465 and 43373).
$ws.Columns("D:E").Insert()

# Copy number formatting (styles) from column F (the old column D, now shifted)
# into the two newly inserted columns D and E, so they inherit the same date /
# number formatting as the rest of the table instead of Excel's bare default.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("F7:F102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows 37 and 79 are section headers with no data columns (only a label in
# column B) - the blanket PasteSpecial above stamped empty, styled cells into
# D37:E37 / D79:E79, so clear those back out to keep those rows label-only.
$ws.Range("D37:E37").Clear()
$ws.Range("D79:E79").Clear()

# Populate the new columns D (period ending 2018-12-31) and E (period ending
# 2018-09-30) with the latest reported financial figures.

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 15100
$ws.Range("E8").Value = 11100
$ws.Range("D9").Value = 1300
$ws.Range("E9").Value = 1000
$ws.Range("D10").Value = 13800
$ws.Range("E10").Value = 10100
$ws.Range("D12").Value = 2900
$ws.Range("E12").Value = 2800
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 700
$ws.Range("E15").Value = 600
$ws.Range("D17").Value = 15000
$ws.Range("E17").Value = 13300
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = -2200
$ws.Range("D20").Value = -100
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 700
$ws.Range("E21").Value = -1600
$ws.Range("D22").Value = 1200
$ws.Range("E22").Value = 1200
$ws.Range("D23").Value = -1200
$ws.Range("E23").Value = -3400
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -1200
$ws.Range("E26").Value = -3500
$ws.Range("D27").Value = -1200
$ws.Range("E27").Value = 34900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 100
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = -1200
$ws.Range("E33").Value = 34900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -1200
$ws.Range("E35").Value = 34900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 13000
$ws.Range("E41").Value = 12600
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 17300
$ws.Range("E43").Value = 14300
$ws.Range("D44").Value = 2400
$ws.Range("E44").Value = 1700
$ws.Range("D45").Value = 2100
$ws.Range("E45").Value = 2800
$ws.Range("D46").Value = 34800
$ws.Range("E46").Value = 31400
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 1400
$ws.Range("E48").Value = 1600
$ws.Range("D49").Value = 16700
$ws.Range("E49").Value = 17200
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1200
$ws.Range("E52").Value = 1000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 54100
$ws.Range("E54").Value = 51200
$ws.Range("D57").Value = 6400
$ws.Range("E57").Value = 5600
$ws.Range("D58").Value = 200
$ws.Range("E58").Value = 200
$ws.Range("D59").Value = 3600
$ws.Range("E59").Value = 1900
$ws.Range("D60").Value = 10200
$ws.Range("E60").Value = 7700
$ws.Range("D61").Value = 38200
$ws.Range("E61").Value = 38000
$ws.Range("D62").Value = 3000
$ws.Range("E62").Value = 2500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 51400
$ws.Range("E66").Value = 48200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 30300
$ws.Range("E70").Value = 30300
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -373400
$ws.Range("E72").Value = -372200
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = -27600
$ws.Range("E76").Value = -27400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -1200
$ws.Range("E81").Value = 34900
$ws.Range("D83").Value = 700
$ws.Range("E83").Value = 600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 700
$ws.Range("E89").Value = -3800
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = -100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = -100
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = -200
$ws.Range("D101").Value = -200
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 500
$ws.Range("E102").Value = -4100
